# Rename worksheets
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("nhap-thanhpham")
$ws1.Name = "nhap-linhkien"

$ws2 = $wb.Worksheets.Item("xuat-thanhpham")
$ws2.Name = "xuat-linhkien"

# New header row, shared across both sheets
$headers = @("Tên Hàng", "Part Number", "Sổ Hợp Đồng", "Sản Phẩm", "Công Ty Nhập", "Ngày Nhập", "Đơn Vị Tính", "Số Lượng", "Đơn Giá", "Thành Tiền")

# --- Sheet1 (nhap-linhkien): only header row, clear old data rows ---
$ws1.Cells.Clear()
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value = $headers[$i]
}
# ColumnWidth is specified in "characters"; 19.2 here round-trips to a
# stored OOXML column width of exactly 20 (matching the original file).
$ws1.Columns("A:J").ColumnWidth = 19.2

# --- Sheet2 (xuat-linhkien): header row + one data row ---
$ws2.Cells.Clear()
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$ws2.Columns("A:J").ColumnWidth = 19.2

$ws2.Range("A2").Value = "món hàng 1"
$ws2.Range("B2").Value = "a"
$ws2.Range("C2").Value = "a"
$ws2.Range("D2").Value = "a"
$ws2.Range("E2").Value = "a"

# "2021-07-02" looks like a date, so a plain .Value assignment would be
# auto-converted to a date serial number. Enter it as a text formula and
# convert it to a static value via copy / paste-special-values so it ends
# up stored as a genuine shared string (matching the original authoring),
# without leaving behind any new number-format/style definitions.
$ws2.Range("F2").Formula = "=""2021-07-02"""
$ws2.Range("F2").Copy()
$ws2.Range("F2").PasteSpecial(-4163)

$ws2.Range("G2").Value = "kg"
$ws2.Range("H2").Value = 1000
$ws2.Range("I2").Value = 20
$ws2.Range("J2").Value = 20000
